# Fix branch revenue report: the "B8" column header used to read
# "Ngày" (Date); it now reads "Thời gian" (Time) instead.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B8").Value = "Thời gian"

# Reflect the author's final cell selection when the workbook was saved.
$ws.Range("E18").Select()
